# Weekly update: a new Orégano price-report row for Vega Central Mapocho de
# Santiago is inserted at row 28 (dated 2022-02-09 / serial 44601), pushing
# all the previously-existing weekly rows (old rows 28-49) down by one
# (new rows 29-50).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28; this shifts rows 28:49 down to 29:50
# and expands the sheet's used range to row 50.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with this week's report.
$ws.Cells.Item(28, 1).Value  = 9
$ws.Cells.Item(28, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(28, 3).Value  = "Metropolitana"
$ws.Cells.Item(28, 4).Value  = 44601
$ws.Cells.Item(28, 5).Value  = 13
$ws.Cells.Item(28, 6).Value  = 100112029
$ws.Cells.Item(28, 7).Value  = "Orégano"
$ws.Cells.Item(28, 8).Value  = "Sin especificar"
$ws.Cells.Item(28, 9).Value  = "Primera"
$ws.Cells.Item(28, 10).Value = 16
$ws.Cells.Item(28, 11).Value = 9000
$ws.Cells.Item(28, 12).Value = 10000
$ws.Cells.Item(28, 13).Value = 9500
$ws.Cells.Item(28, 14).Value = "$/docena de atados"
$ws.Cells.Item(28, 15).Value = "Región Metropolitana"
$ws.Cells.Item(28, 16).Value = 3167
$ws.Cells.Item(28, 17).Value = 3
$ws.Cells.Item(28, 18).Value = "Hortaliza"
